$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its current location
#    (the empty paragraph right before "KNOWN BUGS :"). Deleting the
#    bookmark leaves a clean, empty paragraph behind.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the last bullet point in the "KNOWN BUGS" list
#    ("Snapshots from FaustLive-1.0 are not compatible with
#    FaustLive-2.0") and append a new bullet after it.
# ------------------------------------------------------------------
$lastBullet = $d.Paragraphs.Item($d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Snapshots from FaustLive-1\.0") {
        $lastBullet = $p
        break
    }
}

$insertionPoint = $lastBullet.Range
$insertionPoint.Collapse(0)
# Insert the new bullet paragraph, plus a trailing placeholder character
# ("X") that will carry the new location for the bookmark. Using a real
# character (instead of inserting a bare paragraph mark) avoids leaving
# a stray empty run behind once the character is removed again below.
$insertionPoint.InsertAfter("`rQuit Menu is sometimes disappearing on some systems`rX")

# Recompute paragraph indices after the insert.
$newBulletIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Quit Menu is sometimes disappearing on some systems") {
        $newBulletIndex = $i
        break
    }
}

$placeholderPara = $d.Paragraphs.Item($newBulletIndex + 1)

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark around the placeholder character,
#    then delete the placeholder character so the paragraph collapses
#    back down to just the bookmark start/end markers.
# ------------------------------------------------------------------
$placeholderStart = $placeholderPara.Range.Start
$bmRange = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$charRange = $d.Range($placeholderStart, $placeholderStart + 1)
$charRange.Delete()

# ------------------------------------------------------------------
# 4. The final paragraph should be a plain paragraph (no list style /
#    numbering, it was only inherited from the preceding bulleted
#    paragraph when the paragraph break was inserted).
# ------------------------------------------------------------------
$placeholderPara.Range.ListFormat.RemoveNumbers()
$placeholderPara.Style = "Normal"
